$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.122.75"
$ws.Range("E2").Value = "  +5.13%  "
$ws.Range("D3").Value = "2.264.83"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'230.29"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("D7").Value = "'63.59"
$ws.Range("E7").Value = "  +4.83%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.446"
$ws.Range("E9").Value = "  +11.09%  "
$ws.Range("D10").Value = "'0.104"
$ws.Range("E10").Value = "  +15.22%  "
$ws.Range("D11").Value = "'56.87"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "'26.33"
$ws.Range("E12").Value = "  +19.25%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "2.599.94"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "'15.73"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "'6.06"
$ws.Range("E16").Value = "  +8.69%  "
$ws.Range("D17").Value = "'0.839"
$ws.Range("E17").Value = "  +5.06%  "
$ws.Range("D18").Value = "2.273.36"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").Value = "43.912.51"
$ws.Range("E19").Value = "  +4.79%  "
$ws.Range("D20").Value = "'0.0000101"
$ws.Range("E20").Value = "  +6.98%  "
$ws.Range("D21").Value = "'73.65"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").Value = "'6.03"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "'254.94"
$ws.Range("E23").Value = "  +5.08%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'2.42"
$ws.Range("E25").Value = "  +3.14%  "
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").Value = "'3.36"
$ws.Range("E27").Value = "  +26.91%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'10.15"
$ws.Range("E28").Value = "  +5.39%  "
$ws.Range("D29").Value = "'171.91"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").Value = "'20.84"
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("D31").Value = "'0.138"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("E34").Value = "  +5.20%  "
$ws.Range("D35").Value = "'4.77"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("D36").Value = "'4.87"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("E37").Value = "  +8.80%  "
$ws.Range("E38").Value = "  +7.15%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "'0.0256"
$ws.Range("E40").Value = "  +5.14%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'8.33"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("E43").Value = "  +8.40%  "
$ws.Range("D44").Value = "'0.0965"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").Value = "'97.81"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("D46").Value = "'0.000212"
$ws.Range("E46").Value = "  -8.55%  "
$ws.Range("D47").Value = "'4.38"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "'10.13"
$ws.Range("E49").Value = "  +18.86%  "
$ws.Range("D50").Value = "1.448.59"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "'2.29"
$ws.Range("E51").Value = "  +4.02%  "
